# dataloader modified and environment data only train/infer supported
#
# Column A on the data sheet held Excel serial-date values (e.g. 43178)
# formatted with a custom "YYYY-MM-DD HH:MM:SS" number format. Replace
# each serial date with a plain YYYYMMDD integer (e.g. 20180319) and
# drop the custom date number format so the cell reverts to the
# workbook's default ("Normal") style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($null -eq $serial) {
        continue
    }
    $d = [DateTime]::FromOADate($serial)
    $newVal = [int]$d.ToString("yyyyMMdd")

    # Drop the custom date-time number format (back to the default style)
    $cell.Style = "Normal"
    $cell.Value = $newVal
}

# Best-effort: ask Excel to drop the now-unused custom date number formats
try { $wb.DeleteNumberFormat("YYYY-MM-DD HH:MM:SS") } catch {}
try { $wb.DeleteNumberFormat("yyyy-mm-dd h:mm:ss") } catch {}
